# Update: Session clear logic and precise source attribution prompts
#
# 1. Rewrite titles/bullets on slides 1-3 to the new "公明新聞" weekly-roundup
#    content.
# 2. Rewrite the notes text on slide 2 and slide 3 with precise source
#    attribution.
# 3. Remove the 4th slide ("今後の焦点：政治はどう動くか") entirely. Deleting the
#    slide also removes its notes page (notesSlide3.xml) and its sldId entry
#    (id="259", r:id="rId11") from the presentation's sldIdLst.

$p = $ppt.ActivePresentation

function Set-TitleText($shape, $text) {
    # Title placeholders only ever hold a single paragraph. Assigning
    # TextRange.Text directly against a lone paragraph causes the host to
    # stamp the new run with an explicit <a:rPr lang="en-US"/>, which the
    # original deck never had. Prepending a throwaway paragraph first (so
    # the assignment lands on a multi-paragraph range) avoids that, and
    # then deleting the leading dummy paragraph leaves a single clean
    # <a:r><a:t>...</a:t></a:r> behind.
    $tf = $shape.TextFrame
    $tf.TextRange.Text = "x`r" + $text
    $tf.TextRange.Paragraphs(1).Delete()
}

# --- Slide 1: title -----------------------------------------------------
$s1 = $p.Slides.Item(1)
Set-TitleText $s1.Shapes.Item(1) "公明新聞の今週１週間の気になる話題をまとめてに関する解説"

# --- Slide 2: title + bullets + notes -----------------------------------
$s2 = $p.Slides.Item(2)
Set-TitleText $s2.Shapes.Item(1) "今週の公明新聞：地域組織の「基礎体力」測定"

$body2 = $s2.Shapes.Item(2).TextFrame.TextRange
$body2.Text = "1`r2`r3"
$body2.Paragraphs(1).Text = "・大阪・高石支部：「有権者比1%」を達成（46,827人中484部購読）"
$body2.Paragraphs(2).Text = "・仙台・宮城野支部：支部会に96人参加、替え歌などで結束強化"
$body2.Paragraphs(3).Text = "・宇都宮：企業経営者による購読事例を紹介（支持層拡大のアピール）"

$s2.NotesPage.Shapes.Item(2).TextFrame.TextRange.Text = "出典：公明新聞（2025年12月22日、25日付）。年末の紙面は組織の引き締めと成果の可視化に重点。"

# --- Slide 3: title + bullets + notes -----------------------------------
$s3 = $p.Slides.Item(3)
Set-TitleText $s3.Shapes.Item(1) "国会論戦に見る「公明新聞」の2つの顔"

$body3 = $s3.Shapes.Item(2).TextFrame.TextRange
$body3.Text = "1`r2"
$body3.Paragraphs(1).Text = "・財政面：「事実上赤字で利益になっていない」との指摘も（3/14 福島委員）"
$body3.Paragraphs(2).Text = "・記録面：1993年の記事を引用し、政治改革の「一貫性」を証明（3/28 中川委員）"

$s3.NotesPage.Shapes.Item(2).TextFrame.TextRange.Text = "出典：衆議院 議事録（2025年3月14日、28日）。機関紙は党のアイデンティティ形成の核心。"

# --- Slide 4: delete entirely (and its notes + sldId) --------------------
$p.Slides.Item(4).Delete()
